$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1: Title "A slide" -> split the run after "A" so the trailing space
# becomes its own run: "A" / " " / "slide"
$sh1 = $s.Shapes.Item(1)
$tr1 = $sh1.TextFrame.TextRange
$c1 = $tr1.Characters(2, 1)
$c1.Text = $c1.Text

# Shape 3: table, cell (1,2) "a table" -> split into "a" / " " / "table"
$sh3 = $s.Shapes.Item(3)
$tbl = $sh3.Table
$cell = $tbl.Cell(1, 2)
$tr3 = $cell.Shape.TextFrame.TextRange
$c3 = $tr3.Characters(2, 1)
$c3.Text = $c3.Text

# Shape 7: TextBox "Plus an image" -> split into "Plus" / " " / "an" / " " / "image"
$sh7 = $s.Shapes.Item(7)
$tr7 = $sh7.TextFrame.TextRange
$c7a = $tr7.Characters(8, 1)
$c7a.Text = $c7a.Text
$c7b = $tr7.Characters(5, 1)
$c7b.Text = $c7b.Text
